# feat(translations): update translations 15 november
#
# Appends 11 new translation-key rows (rows 487-497) to Sheet1. Column A
# holds the raw translation key (a few of them end with a literal "."),
# while column B (English copy) and the not-yet-translated language
# columns G:Q simply mirror the English text.
#
# The shared-string table is filled in the same order the original export
# tool used: the English-copy text for every new row first (columns
# B, G:Q), and only afterwards the column-A key text for every new row -
# so the writes below are intentionally grouped by column, not by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Key = "Contribute to keep your language on top."; Copy = "Contribute to keep your language on top" },
    @{ Key = "Validate to keep your language on top.";   Copy = "Validate to keep your language on top" },
    @{ Key = "Please don't use only numerics or email as username"; Copy = "Please don't use only numerics or email as username" },
    @{ Key = "Only 1000 characters allowed"; Copy = "Only 1000 characters allowed" },
    @{ Key = "Contribute to see your language on top."; Copy = "Contribute to see your language on top" },
    @{ Key = "Validate to see your language on top.";   Copy = "Validate to see your language on top" },
    @{ Key = "We feel the text you entered doesn't match the original text, are you sure about your edit"; Copy = "We feel the text you entered doesn't match the original text, are you sure about your edit" },
    @{ Key = "404 Error"; Copy = "404 Error" },
    @{ Key = "Seems this page doesn't exist"; Copy = "Seems this page doesn't exist" },
    @{ Key = "Visit our homepage"; Copy = "Visit our homepage" },
    @{ Key = "Unspecified location"; Copy = "Unspecified location" }
)

$startRow = 487
$langCols = @("G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Pass 1: English copy into column B for every new row (row order).
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("B$r").Value = $rows[$i].Copy
}

# Pass 2: mirror the same English copy across the language columns.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    foreach ($col in $langCols) {
        $ws.Range("$col$r").Value = $rows[$i].Copy
    }
}

# Pass 3: the key itself into column A for every new row (row order).
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $rows[$i].Key
}

# Newly-created cells otherwise pick up a stray wrap-text style; reset every
# touched cell (not the whole row, which would materialise blank cells in
# the untouched C:F columns) back to the workbook's default, un-styled look
# that every other data row in the sheet uses.
$dataCols = @("A", "B", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    foreach ($col in $dataCols) {
        $ws.Range("$col$r").Style = "Normal"
    }
}

# Move the view/selection to where a user would land after typing the last
# new key into column A and pressing Enter.
$null = $ws.Range("A501").Select()
